$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-05 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-06 Friday", 2)
$d.Content.Find.Execute("312÷5=62, 2", $true, $false, $false, $false, $false, $true, 1, $false, "325÷5=65, 0", 2)
$d.Content.Find.Execute("130÷7=18, 4", $true, $false, $false, $false, $false, $true, 1, $false, "377÷9=41, 8", 2)
$d.Content.Find.Execute("608÷3=202, 2", $true, $false, $false, $false, $false, $true, 1, $false, "814÷3=271, 1", 2)
$d.Content.Find.Execute("635÷2=317, 1", $true, $false, $false, $false, $false, $true, 1, $false, "476÷2=238, 0", 2)
$d.Content.Find.Execute("532÷9=59, 1", $true, $false, $false, $false, $false, $true, 1, $false, "410÷2=205, 0", 2)
$d.Content.Find.Execute("438÷3=146, 0", $true, $false, $false, $false, $false, $true, 1, $false, "511÷9=56, 7", 2)
$d.Content.Find.Execute("630÷3=210, 0", $true, $false, $false, $false, $false, $true, 1, $false, "810÷2=405, 0", 2)
$d.Content.Find.Execute("482÷7=68, 6", $true, $false, $false, $false, $false, $true, 1, $false, "225÷5=45, 0", 2)
$d.Content.Find.Execute("559÷5=111, 4", $true, $false, $false, $false, $false, $true, 1, $false, "990÷3=330, 0", 2)
$d.Content.Find.Execute("234÷7=33, 3", $true, $false, $false, $false, $false, $true, 1, $false, "951÷4=237, 3", 2)
$d.Content.Find.Execute("641÷9=71, 2", $true, $false, $false, $false, $false, $true, 1, $false, "175÷5=35, 0", 2)
$d.Content.Find.Execute("883÷7=126, 1", $true, $false, $false, $false, $false, $true, 1, $false, "844÷9=93, 7", 2)
$d.Content.Find.Execute("133÷4=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "529÷5=105, 4", 2)
$d.Content.Find.Execute("497÷4=124, 1", $true, $false, $false, $false, $false, $true, 1, $false, "276÷2=138, 0", 2)
$d.Content.Find.Execute("178÷5=35, 3", $true, $false, $false, $false, $false, $true, 1, $false, "701÷8=87, 5", 2)
$d.Content.Find.Execute("144÷3=48, 0", $true, $false, $false, $false, $false, $true, 1, $false, "286÷9=31, 7", 2)
$d.Content.Find.Execute("281÷2=140, 1", $true, $false, $false, $false, $false, $true, 1, $false, "216÷9=24, 0", 2)
$d.Content.Find.Execute("119÷3=39, 2", $true, $false, $false, $false, $false, $true, 1, $false, "794÷3=264, 2", 2)
$d.Content.Find.Execute("816÷9=90, 6", $true, $false, $false, $false, $false, $true, 1, $false, "425÷9=47, 2", 2)
$d.Content.Find.Execute("751÷2=375, 1", $true, $false, $false, $false, $false, $true, 1, $false, "712÷4=178, 0", 2)
$d.Content.Find.Execute("744÷5=148, 4", $true, $false, $false, $false, $false, $true, 1, $false, "508÷4=127, 0", 2)
$d.Content.Find.Execute("481÷9=53, 4", $true, $false, $false, $false, $false, $true, 1, $false, "638÷7=91, 1", 2)
$d.Content.Find.Execute("582÷6=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "452÷3=150, 2", 2)
$d.Content.Find.Execute("970÷3=323, 1", $true, $false, $false, $false, $false, $true, 1, $false, "454÷9=50, 4", 2)
$d.Content.Find.Execute("312÷2=156, 0", $true, $false, $false, $false, $false, $true, 1, $false, "110÷8=13, 6", 2)
